# Apply the edit described by the diff:
#  - add a new sheet "2024-06-02" (copy of the "2024-06-01" layout/headers)
#  - make the "current" sheet ("current") the active tab/sheet
#  - fix a typo in "current" sheet B2 (Karp_Kuzmin -> Karp_Kuzminnn)
#  - append a new row (row 3) to "current" for the new patient Klara_Novikova
#  - tidy up the duplicate/unused date-time number format on sheet "2024-06-01"
#    (cells that used the now-removed style keep the same look via the remaining style)
#  - bump one timestamp's precision on sheet "2024-06-01"

$wb = $excel.ActiveWorkbook

$sheetCurrent  = $wb.Worksheets.Item("current")
$sheetSettings = $wb.Worksheets.Item("settings")
$sheetJun01    = $wb.Worksheets.Item("2024-06-01")

# ---------------------------------------------------------------------------
# 1. "current" sheet: fix typo, add the new row for 2024-06-02
# ---------------------------------------------------------------------------
$sheetCurrent.Range("B2").Value = "Karp_Kuzminnn"

$sheetCurrent.Range("A3").Value = "2024-06-02"
$sheetCurrent.Range("B3").Value = "Klara_Novikova"
$sheetCurrent.Range("C3").Value = -1
$sheetCurrent.Range("D3").Value = "4"

$sheetCurrent.Range("C17").Select()
$sheetCurrent.Range("F10").Select()

# ---------------------------------------------------------------------------
# 2. Add the new daily sheet "2024-06-02" after "2024-06-01", with the same
#    header row as the other dated sheets.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $sheetJun01)
$newSheet.Name = "2024-06-02"

$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "Время"
$newSheet.Range("C1").Value = "ФИО пациента"
$newSheet.Range("D1").Value = "М\Ж\Р"
$newSheet.Range("E1").Value = "Дата рождения"
$newSheet.Range("F1").Value = "Причина"
$newSheet.Range("G1").Value = "Давление"



# ---------------------------------------------------------------------------
# 3. Bump the timestamp precision for id=7 on "2024-06-01".
# ---------------------------------------------------------------------------
$sheetJun01.Range("B8").Value = 45444.72288207176

# ---------------------------------------------------------------------------
# 4. Make "current" the active sheet/tab (was "2024-06-01").
# ---------------------------------------------------------------------------
$sheetCurrent.Activate()
$sheetCurrent.Select()
